$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TXHD")

# Insert a new column before column D. This shifts the existing D:K
# (quarterly figures) one column to the right, to E:L, and opens up a
# blank column D for the new (most recent) quarter's figures.
$ws.Columns("D").Insert()

# The newly inserted column D inherits formatting from the column to its
# left (column C) by default. Copy the number/cell formatting from the
# neighboring column E (which holds the previous quarter's formatting,
# identical across the whole D:L block in every row) onto the new column D
# so the new quarter's cells are formatted the same way (dates in row
# 7/38/80, plain numbers elsewhere).
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("D7").Value = 43373
$ws.Range("D8").Value = 200
$ws.Range("D9").Value = 100
$ws.Range("D10").Value = 100
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 100
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 400
$ws.Range("D18").Value = -200
$ws.Range("D20").Value = 0
$ws.Range("D21").Value = "NA"
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = -200
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -200
$ws.Range("D27").Value = -200
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("D33").Value = -200
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -200
$ws.Range("D38").Value = 43373
$ws.Range("D41").Value = 0
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 500
$ws.Range("D48").Value = 0
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 0
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 500
$ws.Range("D57").Value = 300
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 100
$ws.Range("D60").Value = 400
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = "NA"
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 400
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -15200
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 100
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43373
$ws.Range("D81").Value = -200
$ws.Range("D83").Value = "NA"
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 0
$ws.Range("D91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 0
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 0
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 0
